$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the pT value in column A from 10.6 to 3.65 for data rows 2-7
$ws.Range("A2:A7").Value = 3.65

# Update the shared string label used in column H (rows 2-7) from
# "AUL-0-PT" to "AUL-0-PT-INT"
$ws.Range("H2:H7").Value = "AUL-0-PT-INT"
